# Correct overlapping column issue in training.py
#
# The report used to have 10 columns (A:J). A new "Date of Last Update"
# column needs to be inserted (now column I), which pushes the old
# "Location of Issue" / "Training Issues" columns out to K / L - but the
# generator script also duplicated the TicketID column into the new J
# column. In addition, a new ticket (257) was appended above the existing
# ticket (253), so the data now spans rows 2:3 instead of just row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Force text storage even for numeric-looking strings (e.g. "257"),
    # then drop back to the default "Normal" style so no stray
    # number-format survives on the cell itself.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "TicketID"
$ws.Range("B1").Value = "TicketType"
$ws.Range("C1").Value = "State"
$ws.Range("D1").Value = "Creation Date"
$ws.Range("E1").Value = "Originator Name"
$ws.Range("F1").Value = "Source"
$ws.Range("G1").Value = "Agent Name"
$ws.Range("H1").Value = "Time Worked (Minutes)"
$ws.Range("I1").Value = "Date of Last Update"
$ws.Range("J1").Value = "TicketID"
$ws.Range("K1").Value = "Location of Issue"
$ws.Range("L1").Value = "Training Issues"

# New header cells (I1, J1, K1, L1) should look like the rest of the bold
# header row. Copy the existing header formatting instead of re-applying
# Bold (which would build a brand-new font record instead of reusing the
# workbook's existing header style).
$ws.Range("A1").Copy()
$ws.Range("I1:L1").PasteSpecial(-4122)

# ---- New row 2: ticket 257 ----
Set-TextCell $ws.Range("A2") "257"
$ws.Range("B2").Value = "Training"
$ws.Range("C2").Value = "In Progress"
$ws.Range("D2").Value = "2017-08-15T11:45:35.400000"
$ws.Range("E2").Value = "Vincent Chov"
$ws.Range("F2").Value = "Website"
$ws.Range("G2").Value = "Vincent Chov"
$ws.Range("H2").Value = "None"
$ws.Range("I2").Value = "2017-08-15T11:46:35.427000"
Set-TextCell $ws.Range("J2") "257"
$ws.Range("K2").Value = "FL - PSI West Office"
$ws.Range("L2").Value = "ANSYS Workbench"

# ---- Row 3 (previously row 2): ticket 253 ----
Set-TextCell $ws.Range("A3") "253"
$ws.Range("B3").Value = "Training"
$ws.Range("C3").Value = "In Progress"
$ws.Range("D3").Value = "2017-08-15T09:31:56.060000"
$ws.Range("E3").Value = "Vincent Chov"
$ws.Range("F3").Value = "Website"
$ws.Range("G3").Value = "Vincent Chov"
$ws.Range("H3").Value = "None"
$ws.Range("I3").Value = "2017-08-15T09:32:25.523000"
Set-TextCell $ws.Range("J3") "253"
$ws.Range("K3").Value = "CT - PSI Hartford Office"
$ws.Range("L3").Value = "Training Content / Topic"
